$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2024-03-06 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-03-07 Thursday", 2) | Out-Null

# Update table cell values
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "90-69=21"
$t.Cell(1,2).Range.Text = "49-46=3"
$t.Cell(1,3).Range.Text = "75-43=32"
$t.Cell(1,4).Range.Text = "92-25=67"
$t.Cell(1,5).Range.Text = "91-44=47"

$t.Cell(2,1).Range.Text = "70-13=57"
$t.Cell(2,2).Range.Text = "13+62=75"
$t.Cell(2,3).Range.Text = "48-41=7"
$t.Cell(2,4).Range.Text = "96-6=90"
$t.Cell(2,5).Range.Text = "67-56=11"

$t.Cell(3,1).Range.Text = "55-28=27"
$t.Cell(3,2).Range.Text = "76+10=86"
$t.Cell(3,3).Range.Text = "39+54=93"
$t.Cell(3,4).Range.Text = "36+61=97"
$t.Cell(3,5).Range.Text = "96-39=57"

$t.Cell(4,1).Range.Text = "17+33=50"
$t.Cell(4,2).Range.Text = "19+74=93"
$t.Cell(4,3).Range.Text = "32+24=56"
$t.Cell(4,4).Range.Text = "61+32=93"
$t.Cell(4,5).Range.Text = "41+36=77"

$t.Cell(5,1).Range.Text = "55-46=9"
$t.Cell(5,2).Range.Text = "9+35=44"
$t.Cell(5,3).Range.Text = "52+7=59"
$t.Cell(5,4).Range.Text = "63-51=12"
$t.Cell(5,5).Range.Text = "28+34=62"

$t.Cell(6,1).Range.Text = "21+35=56"
$t.Cell(6,2).Range.Text = "72-61=11"
$t.Cell(6,3).Range.Text = "7+92=99"
$t.Cell(6,4).Range.Text = "59-6=53"
$t.Cell(6,5).Range.Text = "15+72=87"

$t.Cell(7,1).Range.Text = "0+28=28"
$t.Cell(7,2).Range.Text = "40+42=82"
$t.Cell(7,3).Range.Text = "5+82=87"
$t.Cell(7,4).Range.Text = "79-70=9"
$t.Cell(7,5).Range.Text = "0+26=26"

$t.Cell(8,1).Range.Text = "87-84=3"
$t.Cell(8,2).Range.Text = "36+25=61"
$t.Cell(8,3).Range.Text = "57+28=85"
$t.Cell(8,4).Range.Text = "66+30=96"
$t.Cell(8,5).Range.Text = "60-23=37"

$t.Cell(9,1).Range.Text = "94-34=60"
$t.Cell(9,2).Range.Text = "8+21=29"
$t.Cell(9,3).Range.Text = "16+70=86"
$t.Cell(9,4).Range.Text = "5+74=79"
$t.Cell(9,5).Range.Text = "45-31=14"

$t.Cell(10,1).Range.Text = "40+26=66"
$t.Cell(10,2).Range.Text = "89+6=95"
$t.Cell(10,3).Range.Text = "88-14=74"
$t.Cell(10,4).Range.Text = "25+29=54"
$t.Cell(10,5).Range.Text = "89-69=20"

$t.Cell(11,1).Range.Text = "51-4=47"
$t.Cell(11,2).Range.Text = "27+64=91"
$t.Cell(11,3).Range.Text = "94+4=98"
$t.Cell(11,4).Range.Text = "3+79=82"
$t.Cell(11,5).Range.Text = "53+35=88"

$t.Cell(12,1).Range.Text = "48+6=54"
$t.Cell(12,2).Range.Text = "97-49=48"
$t.Cell(12,3).Range.Text = "87-9=78"
$t.Cell(12,4).Range.Text = "86-37=49"
$t.Cell(12,5).Range.Text = "85-3=82"

$t.Cell(13,1).Range.Text = "30+0=30"
$t.Cell(13,2).Range.Text = "52+47=99"
$t.Cell(13,3).Range.Text = "39+6=45"
$t.Cell(13,4).Range.Text = "54-38=16"
$t.Cell(13,5).Range.Text = "11+38=49"

$t.Cell(14,1).Range.Text = "16+26=42"
$t.Cell(14,2).Range.Text = "24+0=24"
$t.Cell(14,3).Range.Text = "37+16=53"
$t.Cell(14,4).Range.Text = "21+12=33"
$t.Cell(14,5).Range.Text = "54+36=90"

$t.Cell(15,1).Range.Text = "85-84=1"
$t.Cell(15,2).Range.Text = "45-9=36"
$t.Cell(15,3).Range.Text = "25+46=71"
$t.Cell(15,4).Range.Text = "23+67=90"
$t.Cell(15,5).Range.Text = "42-29=13"

$t.Cell(16,1).Range.Text = "49-26=23"
$t.Cell(16,2).Range.Text = "14+84=98"
$t.Cell(16,3).Range.Text = "32-2=30"
$t.Cell(16,4).Range.Text = "90-75=15"
$t.Cell(16,5).Range.Text = "45-3=42"

$t.Cell(17,1).Range.Text = "16+68=84"
$t.Cell(17,2).Range.Text = "20+5=25"
$t.Cell(17,3).Range.Text = "63-16=47"
$t.Cell(17,4).Range.Text = "84-24=60"
$t.Cell(17,5).Range.Text = "3+14=17"

$t.Cell(18,1).Range.Text = "63+9=72"
$t.Cell(18,2).Range.Text = "63-26=37"
$t.Cell(18,3).Range.Text = "3+23=26"
$t.Cell(18,4).Range.Text = "93-13=80"
$t.Cell(18,5).Range.Text = "28+60=88"

$t.Cell(19,1).Range.Text = "46+13=59"
$t.Cell(19,2).Range.Text = "76+6=82"
$t.Cell(19,3).Range.Text = "78-0=78"
$t.Cell(19,4).Range.Text = "5+6=11"
$t.Cell(19,5).Range.Text = "68+2=70"

$t.Cell(20,1).Range.Text = "26+14=40"
$t.Cell(20,2).Range.Text = "90+0=90"
$t.Cell(20,3).Range.Text = "39+12=51"
$t.Cell(20,4).Range.Text = "61-33=28"
$t.Cell(20,5).Range.Text = "88-78=10"
